$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 - this shifts existing rows 32:63 down to 33:64
# and expands the used range to A1:R64 (matching the diff's dimension change).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with a new weekly price entry for Albahaca
# (same market/category/quality metadata as the rest of the block; only the
# date and price figures differ).
$ws.Cells.Item(32, 1).Value = 1
$ws.Cells.Item(32, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value = 45096
$ws.Cells.Item(32, 5).Value = 15
$ws.Cells.Item(32, 6).Value = 100112052
$ws.Cells.Item(32, 7).Value = "Albahaca"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 250
$ws.Cells.Item(32, 11).Value = 900
$ws.Cells.Item(32, 12).Value = 1000
$ws.Cells.Item(32, 13).Value = 950
$ws.Cells.Item(32, 14).Value = "$/paquete"
$ws.Cells.Item(32, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(32, 16).Value = 950
$ws.Cells.Item(32, 17).Value = 1
$ws.Cells.Item(32, 18).Value = "Hortaliza"
